# Add two new rows of daily error-count data to Sheet1, mirroring the
# manual "Add files via upload" edit: rows for 1/27/2026 and 1/28/2026.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data rows (date serials 46049 and 46050 correspond to 2026-01-27 and 2026-01-28)
$ws.Range("A82").Value = 46049
$ws.Range("B82").Value = 936
$ws.Range("C82").Value = 916
$ws.Range("D82").Value = 20

$ws.Range("A83").Value = 46050
$ws.Range("B83").Value = 3044
$ws.Range("C83").Value = 3040
$ws.Range("D83").Value = 4

# Match date number formatting used by the rest of column A (style index 3 => numFmtId 15)
$ws.Range("A82:A83").NumberFormat = $ws.Range("A81").NumberFormat

# Update the visible viewport/selection as Excel would after scrolling to the new rows
$ws.Range("V92").Select()
$excel.ActiveWindow.ScrollRow = 65
